# Natmi following Dr Hou advice:
# Re-run the sending/target cluster cross-product for Lama2-Dag1 to include the "ECs" cluster
# (previously only FAPs and sCs were senders/targets; now ECs, FAPs, sCs all appear on both sides).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row data for rows 2-10 (A:T), in column order A,B,C,D,E,F,G,H,I,J,K,L,M,N,O,P,Q,R,S,T
$rowData = @(
  @("ECs", "Lama2", "Dag1", "ECs", 2, 0.6666666666666666, 2.452389333333333, 7.357168, 0.007993767302975028, 0.007993767302975028, 3, 1, 14.792724, 44.378172, 0.0761423615956231, 0.0761423615956231, 36.27751854854399, 326.497666936896, 0.0006086643204943934, 0.0006086643204943934),
  @("ECs", "Lama2", "Dag1", "FAPs", 2, 0.6666666666666666, 2.452389333333333, 7.357168, 0.007993767302975028, 0.007993767302975028, 3, 1, 44.28072733333334, 132.842182, 0.2279255093471533, 0.2279255093471533, 108.5935833845084, 977.342250460576, 0.001821983484133203, 0.001821983484133203),
  @("ECs", "Lama2", "Dag1", "sCs", 2, 0.6666666666666666, 2.452389333333333, 7.357168, 0.007993767302975028, 0.007993767302975028, 3, 1, 135.203738, 405.611214, 0.6959321290572236, 0.6959321290572236, 331.5722048979947, 2984.149844081952, 0.005563119498347431, 0.005563119498347431),
  @("FAPs", "Lama2", "Dag1", "ECs", 3, 1, 243.5672963333334, 730.701889, 0.7939278902575405, 0.7939278902575405, 3, 1, 14.792724, 44.378172, 0.0761423615956231, 0.0761423615956231, 3603.023790085212, 32427.21411076691, 0.06045154450083982, 0.06045154450083982),
  @("FAPs", "Lama2", "Dag1", "FAPs", 3, 1, 243.5672963333334, 730.701889, 0.7939278902575405, 0.7939278902575405, 3, 1, 44.28072733333334, 132.842182, 0.2279255093471533, 0.2279255093471533, 10785.33703625354, 97068.03332628182, 0.1809564187718607, 0.1809564187718607),
  @("FAPs", "Lama2", "Dag1", "sCs", 3, 1, 243.5672963333334, 730.701889, 0.7939278902575405, 0.7939278902575405, 3, 1, 135.203738, 405.611214, 0.6959321290572236, 0.6959321290572236, 32931.20891882037, 296380.8802693833, 0.5525199269848399, 0.5525199269848399),
  @("sCs", "Lama2", "Dag1", "ECs", 3, 1, 60.76799533333334, 182.303986, 0.1980783424394845, 0.1980783424394845, 3, 1, 14.792724, 44.378172, 0.0761423615956231, 0.0761423615956231, 898.9241829992881, 8090.317646993592, 0.01508215277428888, 0.01508215277428888),
  @("sCs", "Lama2", "Dag1", "FAPs", 3, 1, 60.76799533333334, 182.303986, 0.1980783424394845, 0.1980783424394845, 3, 1, 44.28072733333334, 132.842182, 0.2279255093471533, 0.2279255093471533, 2690.851031948606, 24217.65928753746, 0.04514710709115935, 0.04514710709115935),
  @("sCs", "Lama2", "Dag1", "sCs", 3, 1, 60.76799533333334, 182.303986, 0.1980783424394845, 0.1980783424394845, 3, 1, 135.203738, 405.611214, 0.6959321290572236, 0.6959321290572236, 8216.060119833224, 73944.54107849901, 0.1378490825740362, 0.1378490825740362)
)

$startRow = 2
for ($i = 0; $i -lt $rowData.Length; $i++) {
  $values = $rowData[$i]
  $r = $startRow + $i
  for ($c = 0; $c -lt $values.Length; $c++) {
    $ws.Cells.Item($r, $c + 1).Value = $values[$c]
  }
}
